$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "69.047.40"
$ws.Range("E2").Value = "  +0.52%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.765.05"
$ws.Range("E3").Value = "  -1.42%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.13%  "

# Row 5 - BNB
Set-TextValue "D5" "631.66"
$ws.Range("E5").Value = "  +3.00%  "

# Row 6 - Solana
Set-TextValue "D6" "165.98"
$ws.Range("E6").Value = "  +0.71%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.762.35"
$ws.Range("E7").Value = "  -1.49%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.03%  "

# Row 9 - XRP
Set-TextValue "D9" "0.520"
$ws.Range("E9").Value = "  +0.45%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -1.52%  "

# Row 11 - Cardano
Set-TextValue "D11" "0.458"
$ws.Range("E11").Value = "  +1.49%  "

# Row 12 - Toncoin
$ws.Range("E12").Value = "  +1.02%  "

# Row 13 - ShibaInu
Set-TextValue "D13" "0.0000239"
$ws.Range("E13").Value = "  -3.52%  "

# Row 14 - Avalanche
Set-TextValue "D14" "34.84"
$ws.Range("E14").Value = "  -1.75%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "4.404.08"
$ws.Range("E15").Value = "  -1.29%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.762.74"
$ws.Range("E16").Value = "  -1.67%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "69.081.73"
$ws.Range("E17").Value = "  +0.61%  "

# Row 18 - Chainlink
Set-TextValue "D18" "17.61"
$ws.Range("E18").Value = "  -2.86%  "

# Row 19 - TRON
$ws.Range("E19").Value = "  +0.15%  "

# Row 20 - Polkadot
$ws.Range("E20").Value = "  -1.55%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "461.75"
$ws.Range("E21").Value = "  -0.69%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -1.43%  "

# Row 23 - Polygon
Set-TextValue "D23" "0.706"
$ws.Range("E23").Value = "  +0.65%  "

# Row 24 - PEPE
Set-TextValue "D24" "0.0000144"
$ws.Range("E24").Value = "  -3.53%  "

# Row 25 - Litecoin
Set-TextValue "D25" "82.12"
$ws.Range("E25").Value = "  -2.01%  "

# Row 26 - InternetComputer(DFINITY)
Set-TextValue "D26" "12.12"
$ws.Range("E26").Value = "  +0.75%  "

# Row 27 - Fetch.AI
Set-TextValue "D27" "2.13"
$ws.Range("E27").Value = "  +0.40%  "

# Row 28 - RenderToken
Set-TextValue "D28" "10.07"
$ws.Range("E28").Value = "  +0.71%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  -0.04%  "

# Row 30 - WrappedeETH
Set-TextValue "D30" "3.917.93"
$ws.Range("E30").Value = "  -1.21%  "

# Row 31 - ImmutableX
$ws.Range("E31").Value = "  +3.11%  "

# Row 32 - PancakeSwap
Set-TextValue "D32" "2.68"
$ws.Range("E32").Value = "  +1.85%  "

# Row 33 - NEARProtocol
Set-TextValue "D33" "7.05"

# Row 34 - now EthereumClassic (was Kaspa)
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D34" "28.40"
$ws.Range("E34").Value = "  -2.17%  "

# Row 35 - now Kaspa (was EthereumClassic)
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D35" "0.175"
$ws.Range("E35").Value = "  +19.06%  "

# Row 36 - Binance-PegBSC-USD
$ws.Range("E36").Value = "  -0.04%  "

# Row 37 - RenzoRestakedETH
Set-TextValue "D37" "3.722.80"
$ws.Range("E37").Value = "  -1.18%  "

# Row 38 - Aptos
Set-TextValue "D38" "8.93"
$ws.Range("E38").Value = "  -1.48%  "

# Row 39 - Hedera
Set-TextValue "D39" "0.100"
$ws.Range("E39").Value = "  -1.03%  "

# Row 40 - dogwifhat
Set-TextValue "D40" "3.28"
$ws.Range("E40").Value = "  +4.13%  "

# Row 41 - Filecoin
Set-TextValue "D41" "5.79"
$ws.Range("E41").Value = "  -1.88%  "

# Row 42 - Mantle
Set-TextValue "D42" "0.964"
$ws.Range("E42").Value = "  -1.65%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  +0.09%  "

# Row 44 - USDe
$ws.Range("E44").Value = "  -0.03%  "

# Row 45 - Monero
Set-TextValue "D45" "157.34"
$ws.Range("E45").Value = "  +2.39%  "

# Row 46 - Stacks
$ws.Range("E46").Value = "  +5.00%  "

# Row 47 - ONDO
$ws.Range("E47").Value = "  +2.82%  "

# Row 48 - OKB
Set-TextValue "D48" "46.97"
$ws.Range("E48").Value = "  +0.87%  "

# Row 49 - Arweave
Set-TextValue "D49" "43.03"
$ws.Range("E49").Value = "  +0.91%  "

# Row 50 - TheGraph
Set-TextValue "D50" "0.294"
$ws.Range("E50").Value = "  -1.32%  "

# Row 51 - Cosmos
Set-TextValue "D51" "8.35"
$ws.Range("E51").Value = "  -0.29%  "
